$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (DataFetchFlag): add rows 13, 12, 14 (this order matches the
# order in which the new shared strings were first introduced upstream) ---

$ws1.Range("A13").Value = "FLD_Transmittals_ActionRequired_New_RequestForInformation"
$ws1.Range("B13").Value = "XL"

$ws1.Range("A12").Value = "FLD_Transmittals_ActionRequired_New_IssuedForReview"
$ws1.Range("B12").Value = "XL"

$ws1.Range("A14").Value = "FLD_Transmittals_ActionRequired_New_IssuedForApproval"
$ws1.Range("B14").Value = "XL"

# Re-create the list data validation so it covers the expanded range B2:B14
$ws1.Range("B2:B14").Validation.Delete()
$ws1.Range("B2:B14").Validation.Add(3, 1, 1, '"XL,DB"')

# --- Sheet2 (DataFetchXL): add rows 13, 12, 14 with hyperlinks ---

$ws2.Range("A13").Value = "FLD_Transmittals_ActionRequired_New_RequestForInformation"
$ws2.Hyperlinks.Add($ws2.Range("B13"), 'file:///\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-RequestForInformation.xlsx', "", "", '\\src\\com\\proj\\suiteTRANSMITTALS\\testdata\\TransmittalsTestData-RequestForInformation.xlsx')
$ws2.Range("C13").Value = "Transmittals_New"

$ws2.Range("A12").Value = "FLD_Transmittals_ActionRequired_New_IssuedForReview"
$ws2.Hyperlinks.Add($ws2.Range("B12"), 'file:///\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-IssuedForReview.xlsx', "", "", '\\src\\com\\proj\\suiteTRANSMITTALS\\testdata\\TransmittalsTestData-IssuedForReview.xlsx')
$ws2.Range("C12").Value = "Transmittals_New"

$ws2.Range("A14").Value = "FLD_Transmittals_ActionRequired_New_IssuedForApproval"
$ws2.Hyperlinks.Add($ws2.Range("B14"), 'file:///\\src\com\proj\suiteTRANSMITTALS\testdata\TransmittalsTestData-IssuedForApproval.xlsx', "", "", '\\src\\com\\proj\\suiteTRANSMITTALS\\testdata\\TransmittalsTestData-IssuedForApproval.xlsx')
$ws2.Range("C14").Value = "Transmittals_New"
